$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": add a new daily data row (5/11/2020) above the footnote row,
# pushing the footnote from row 34 down to row 35.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows("34").Insert()
$wsAll.Range("A34").Value = 43962
$wsAll.Range("B34").Value = 278
$wsAll.Range("C34").Value = 274
$wsAll.Range("D34").Value = 80
$wsAll.Range("E34").Value = 69
$wsAll.Range("F34").Value = 11
$wsAll.Range("G34").Value = 8
$wsAll.Range("H34").Value = 186

# ---------------------------------------------------------------------------
# Sheet "kobe": correct the previous day's row (5/10/2020) and add a new
# daily data row (5/11/2020) above the footnote row, pushing the footnote
# from row 89 down to row 90.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Range("D88").Value = 1
$wsKobe.Range("E88").Value = 278

$wsKobe.Rows("89").Insert()
$wsKobe.Range("A89").Value = 43962
$wsKobe.Range("B89").Value = 0
$wsKobe.Range("C89").Value = 2600
$wsKobe.Range("D89").Value = 0
$wsKobe.Range("E89").Value = 278
$wsKobe.Range("F89").Value = 75
$wsKobe.Range("G89").Value = 65
$wsKobe.Range("H89").Value = 10
$wsKobe.Range("I89").Value = 8
$wsKobe.Range("J89").Value = 177

# ---------------------------------------------------------------------------
# Sheet "other": add a new daily data row (5/11/2020) above the footnote
# row, pushing the footnote from row 64 down to row 65.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows("64").Insert()
$wsOther.Range("A64").Value = 43962
$wsOther.Range("B64").Value = 0
$wsOther.Range("C64").Value = 14
$wsOther.Range("D64").Value = 5
$wsOther.Range("E64").Value = 4
$wsOther.Range("F64").Value = 1
$wsOther.Range("G64").Value = 0
$wsOther.Range("H64").Value = 9

Write-Output "edit applied"
